# "fixed bugs in performane view"
# The marks sheet listed the wrong student names and one stale "0/30" mark.
# Correct the roster: update the three student names and make sure every
# mark reads "Not Completed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MarksSheet")

$ws.Range("A2").Value = "W.A.T.N.Jayathilake"
$ws.Range("B2").Value = "Not Completed"

$ws.Range("A3").Value = "D.W.S.N.Sewwandi"
$ws.Range("B3").Value = "Not Completed"

$ws.Range("A4").Value = "L.R.M.U.BANDARA"
$ws.Range("B4").Value = "Not Completed"
